$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the RPA dataset (IPO underwriting table) to the 2024-05-21 push.
# Net effect vs. the previous snapshot:
#   - row 4/5 (KB 민테크 / 제일엠앤에스) swap order
#   - row 11/12 (신한 신한제12호스팩 / HD현대마린솔루션) swap order
#   - row 18/19 (하나 HD현대마린솔루션 / 하나33호스팩) swap order
#   - the 하나32호스팩 row (row 17) drops out of the dataset entirely
#
# Rows are moved with Copy+Insert (and the now-duplicated row removed with
# Delete) rather than via Value/Value2 round-trips, so date-looking text
# such as "2024-04-23" is relocated verbatim instead of being re-parsed by
# Excel into a date serial number.

# Swap row 4 and row 5
$ws.Rows.Item(5).Copy()
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(6).Delete()

# Swap row 11 and row 12
$ws.Rows.Item(12).Copy()
$ws.Rows.Item(11).Insert()
$ws.Rows.Item(13).Delete()

# Swap row 18 and row 19
$ws.Rows.Item(19).Copy()
$ws.Rows.Item(18).Insert()
$ws.Rows.Item(20).Delete()

# Remove the 하나32호스팩 row (old row 17); everything below shifts up.
$ws.Rows.Item(17).Delete()
